$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Populate Sheet2 with the "Names" lookup table (Name / Item / Value)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws2.Range("A1").Value = "Name"
$ws2.Range("B1").Value = "Item"
$ws2.Range("C1").Value = "Value"

$ws2.Range("A2").Value = "Will"
$ws2.Range("B2").Value = "ball"
$ws2.Range("C2").Value = 2

$ws2.Range("A3").Value = "Will"
$ws2.Range("B3").Value = "bat"
$ws2.Range("C3").Value = 3

$ws2.Range("A4").Value = "Mark"
$ws2.Range("B4").Value = "glove"
$ws2.Range("C4").Value = 2

$ws2.Range("A5").Value = "Grove"
$ws2.Range("B5").Value = "boot"
$ws2.Range("C5").Value = 3

$ws2.Range("A6").Value = "Gary"
$ws2.Range("B6").Value = "hat"
$ws2.Range("C6").Value = 2

$ws2.Range("A7").Value = "Gary"
$ws2.Range("B7").Value = "boots"
$ws2.Range("C7").Value = 3

$ws2.Range("A8").Value = "Scary"

# Turn the A1:C8 range into a table (Table13), matching TableStyleDark3,
# no auto-filter arrows shown.
$tbl2 = $ws2.ListObjects.Add(1, $ws2.Range("A1:C8"), $null, 1)
$tbl2.Name = "Table13"
$tbl2.TableStyle = "TableStyleDark3"
$tbl2.ShowAutoFilter = $false

# ---------------------------------------------------------------------------
# 2. Re-point the Sheet1 "unique name" array formulas (J2:J22) at the new
#    Sheet2 table instead of the old $L$10:$L$19 / $L$10:$L$14 ranges, and
#    extend the series down through row 22.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

for ($r = 2; $r -le 22; $r++) {
    $prev = $r - 1
    $formula = "=IFERROR(INDEX(Sheet2!`$A`$2:`$A`$8, MATCH(0, COUNTIF(`$J`$1:J$prev, Sheet2!`$A`$2:`$A`$8), 0)), `"`")"
    $ws1.Range("J$r").FormulaArray = $formula
}

# ---------------------------------------------------------------------------
# 3. Update the Sheet1 selection and restore Sheet1 as the active sheet
#    (selecting on Sheet2 above would otherwise steal "tabSelected").
# ---------------------------------------------------------------------------
$ws1.Range("J2:J22").Select()
$ws1.Activate()

# Leave the Sheet2 selection on B8, as recorded after populating the table.
$ws2.Range("B8").Select()
$ws1.Activate()
